$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$p.Range.Delete()
